$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Hour-conversion block (J2:P8), mirrors B2:H8 * 24 ---
# Row 2: type in J2, then fill right into K2:P2
$ws.Range("J2").Formula = "=B2*24"
$ws.Range("K2:P2").Formula = "=C2*24"

# Rows 3-8: fill the J2:P2 row down column-by-column (J3:J8, K3:K8, ... P3:P8)
$ws.Range("J3:J8").Formula = "=B3*24"
$ws.Range("K3:K8").Formula = "=C3*24"
$ws.Range("L3:L8").Formula = "=D3*24"
$ws.Range("M3:M8").Formula = "=E3*24"
$ws.Range("N3:N8").Formula = "=F3*24"
$ws.Range("O3:O8").Formula = "=G3*24"
$ws.Range("P3:P8").Formula = "=H3*24"

# Strip the inherited time-format style picked up from the source B:H columns
# so the new J:P cells keep the default (General) style, as in the target file.
$ws.Range("J2:P8").Style = "Normal"

# --- Rounded-to-quarter-hour block (J10:P16), mirrors J2:P8 rounded ---
$ws.Range("J10").Formula = "=ROUND(J2*4, 0)/4"
$ws.Range("K10:P10").Formula = "=ROUND(K2*4, 0)/4"

$ws.Range("J11:P11").Formula = "=ROUND(J3*4, 0)/4"
$ws.Range("J12:P12").Formula = "=ROUND(J4*4, 0)/4"
$ws.Range("J13:P13").Formula = "=ROUND(J5*4, 0)/4"
$ws.Range("J14:P14").Formula = "=ROUND(J6*4, 0)/4"
$ws.Range("J15:P15").Formula = "=ROUND(J7*4, 0)/4"
$ws.Range("J16:P16").Formula = "=ROUND(J8*4, 0)/4"

$ws.Range("J10:P16").Style = "Normal"

# --- View state: scroll so column H is leftmost, select J10 ---
$ws.Application.ActiveWindow.ScrollColumn = 8
$ws.Range("J10").Select()
